$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (dSF) updates per the diff
$ws.Range("F2").Value  = -11
$ws.Range("F5").Value  = -11
$ws.Range("F6").Value  = -2
$ws.Range("F7").Value  = -6
$ws.Range("F8").Value  = -1
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = -5
$ws.Range("F14").Value = -1
$ws.Range("F15").Value = -3
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 12
$ws.Range("F19").Value = -6
